# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" worksheet (fund-holding detail) right before
#    the "总计" (totals) summary sheet.
# 2. Insert a new top data row into "总计" summarising the 2022-Q1 sheet
#    (date / holding count / holding market value), shifting the existing
#    rows down and renumbering the index column.

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("2021-Q4")
$totalBeforeAdd = $wb.Worksheets.Item("总计")

# --- 1. Add the "2022-Q1" sheet, placed immediately before "总计" ------------
$newSheet = $wb.Worksheets.Add($totalBeforeAdd)
$newSheet.Name = "2022-Q1"

# NOTE: sheet references in this runtime are positional anchors, not stable
# object handles - now that "总计" has been pushed one slot further back by
# the insert above, we must re-resolve it (by name) to get a handle that
# actually points at the totals sheet instead of the freshly added one.
$total = $wb.Worksheets.Item("总计")

# Copy header formatting (bold / bordered) and the index-column formatting
# from the "2021-Q4" sheet, which already uses the same 7-column layout.
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$newSheet.Range("A2:A8").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Fund code (B) and the numeric-looking percentage/weight columns (D:G) are
# stored as text in every quarter sheet - force text so Excel doesn't coerce
# "012150" -> 12150 or "0.3655" -> a number.
$newSheet.Range("B2:B8").NumberFormat = "@"
$newSheet.Range("D2:G8").NumberFormat = "@"

$data = @(
    @("012150", "诺德价值发现一年持有期混合型证券投资基金", "9.42", "91.48", "3.88", "0.3655", 9),
    @("009994", "嘉实创新先锋混合A", "11.99", "91.88", "2.57", "0.3081", 10),
    @("012036", "诺德兴远优选一年持有期混合型证券投资基金", "2.75", "52.19", "2.87", "0.0789", 6),
    @("009995", "嘉实创新先锋混合C", "2.14", "91.88", "2.57", "0.0550", 10),
    @("013441", "西藏东财创新医疗六个月定开混合", "0.58", "81.46", "5.29", "0.0307", 8),
    @("011149", "创金合信ESG责任投资股票A", "0.16", "87.53", "3.59", "0.0057", 9),
    @("011150", "创金合信ESG责任投资股票C", "0.08", "87.53", "3.59", "0.0029", 9)
)

$r = 2
foreach ($row in $data) {
    $newSheet.Cells.Item($r, 1).Value = $r - 2
    $newSheet.Cells.Item($r, 2).Value = $row[0]
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    $newSheet.Cells.Item($r, 4).Value = $row[2]
    $newSheet.Cells.Item($r, 5).Value = $row[3]
    $newSheet.Cells.Item($r, 6).Value = $row[4]
    $newSheet.Cells.Item($r, 7).Value = $row[5]
    $newSheet.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# --- 2. Add the 2022-Q1 summary row at the top of "总计" --------------------
$total.Range("A2:D2").Insert(-4121)
$total.Range("B2:D2").ClearFormats()
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 7
$total.Cells.Item(2, 4).Value = 0.85

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(5, 1).Value = 3
